$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.753.84"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").Value = "3.727.97"
$ws.Range("E3").Value = "  -2.35%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.36%  "

$ws.Range("D7").Value = "3.727.75"
$ws.Range("E7").Value = "  -2.33%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -1.82%  "

$ws.Range("E10").Value = "  -3.88%  "

$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("E12").Value = "  -2.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000262"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.72%  "

$ws.Range("D15").Value = "4.351.37"
$ws.Range("E15").Value = "  -2.43%  "

$ws.Range("D16").Value = "3.716.46"
$ws.Range("E16").Value = "  -3.02%  "

$ws.Range("D17").Value = "67.672.83"
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.71%  "

$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "466.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.701"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.09%  "

$ws.Range("E25").Value = "  -4.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000135"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -11.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.57%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("D30").Value = "3.871.26"
$ws.Range("E30").Value = "  -2.38%  "

$ws.Range("E31").Value = "  -5.35%  "

$ws.Range("E32").Value = "  -2.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.61%  "

$ws.Range("D37").Value = "3.678.47"
$ws.Range("E37").Value = "  -2.79%  "

$ws.Range("E38").Value = "  -5.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.138"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.989"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.89%  "

$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.306"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "390.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.34%  "
